# "Generate Report for Archive"
#
# The localization status report moved from "Ready for handoff" to
# "In Translation" for every tracked file. That shared string shows up
# on the Overview sheet (one column per locale: zh-cn, de-de) and on
# each locale's own detail sheet (its "Status" column).
#
# Updating the cell text also makes the "Status" column narrower
# (Excel auto-shrinks a previously auto-fit column when its longest
# entry gets shorter), so the ColumnWidth is tightened to match.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn detail sheet: Status column (col C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de detail sheet: Status column (col C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
